# Add "Station Area Access Plan" projects (bangalore CMP 2020) to sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$stations = @(
    @{ Id = 4001; Location = "Baiyapannahalli Station";       Coord = "12.990866321508554, 77.65294454951227" },
    @{ Id = 4002; Location = "Yellachanahalli Station";        Coord = "12.895860594419416, 77.57004255678578" },
    @{ Id = 4003; Location = "Nayandanahalli Station";         Coord = "12.946429708785299, 77.52976417747144" },
    @{ Id = 4004; Location = "Nagasandra Station";             Coord = "13.048204340142766, 77.50014317909618" },
    @{ Id = 4005; Location = "Sampige Mantri Square Station";  Coord = "12.99057903833429, 77.57085307602627" },
    @{ Id = 4006; Location = "Sriramapura Station";            Coord = "12.9966337716117, 77.56351455573149" },
    @{ Id = 4007; Location = "Rajajinagar Station";            Coord = "13.000375243187055, 77.5498452946652" },
    @{ Id = 4008; Location = "Mahakavi Kuvempu Road Station";  Coord = "12.99852857334125, 77.55704994598644" },
    @{ Id = 4009; Location = "Mahalakshmi Station";            Coord = "13.008238795034371, 77.54894763067821" },
    @{ Id = 4010; Location = "Sandal Soap Factory Station";    Coord = "13.014845877265085, 77.55412978470413" },
    @{ Id = 4011; Location = "Yeshwantpur Station";            Coord = "13.023334344199698, 77.54960785891855" },
    @{ Id = 4012; Location = "Peenya Station";                 Coord = "13.032955134958728, 77.53343784692174" },
    @{ Id = 4013; Location = "Peenya Industry Station";        Coord = "13.036495649964248, 77.52539754228509" },
    @{ Id = 4014; Location = "Dasarahalli Station";            Coord = "13.043586550497832, 77.51250107757342" }
)

$startRow = 266

# Fill id / type / location columns first (matches the order the shared
# strings were originally authored in).
for ($i = 0; $i -lt $stations.Count; $i++) {
    $row = $startRow + $i
    $s = $stations[$i]
    $ws.Cells.Item($row, 1).Value = $s.Id
    $ws.Cells.Item($row, 2).Value = "Station Area Access Plan"
    $ws.Cells.Item($row, 3).Value = $s.Location
}

# Then fill the lat/long column. Row 278's coordinate was entered before
# row 277's in the source data, so replicate that write order.
$fOrder = @(0,1,2,3,4,5,6,7,8,9,10,12,11,13)
foreach ($i in $fOrder) {
    $row = $startRow + $i
    $s = $stations[$i]
    $ws.Cells.Item($row, 6).Value = $s.Coord
}

# Column B widened to fit the new "Station Area Access Plan" / station
# names (23.1666... character width serializes to the 24 stored width).
$ws.Columns.Item(2).ColumnWidth = 23.1666666666667

# Park the view near the bottom of the newly added rows.
$ws.Range("F280").Select() | Out-Null
